# Paper_details.xlsx - add notes/details for existing "Murdoch et al" paper (row 10)
# and add a brand new paper entry "Jansen et al" (row 11)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 10: fill in the previously empty Subject (D) and Details (E) cells ---
$ws.Range("D10").Value = "Modelling best strategies for conserving Argentine grasslands "
$ws.Range("D10").WrapText = $true

$ws.Range("E10").Value = "Modelled three different strategies for selecting 10% of a huge area of Argentine grasslands. Minimise cost, maximise gain, and then return on investment. The ROI strategy was by far the best. The maximise gain was by far the worst as the costs were extremely high"
$ws.Range("E10").WrapText = $true

# --- Row 11: brand new paper entry ---
$ws.Range("A11").Value = "Jansen et al"

$ws.Range("B11").Value = "Determinants of income-earning strategies and adoption of conservation practices in hillside communities in rural Honduras"
$ws.Range("B11").WrapText = $true

$ws.Range("C11").Value = 2006

$ws.Range("E11").Value = "State that population increases and insecure tenure have caused forest loss. The results regarding the influence of population density on conservation practices suggest a U-type relationship. That is, up to a certain population density the four conservation practices considered here are less common in communities with higher population densities. However, after a certain point population density has a positive influence on the adoption of conservation practices."
$ws.Range("E11").WrapText = $true

$ws.Range("D11").Value = "Modelled people's willingness to invest in conservation-focused livelihoods"
$ws.Range("D11").WrapText = $true

$ws.Range("A11:E11").RowHeight = 60.6

# --- Update view/selection to match where the author ended up editing ---
$win = $excel.ActiveWindow
$win.ScrollRow = 5
$win.ScrollColumn = 1
$ws.Range("D11").Select()
